$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F574/G574 timestamps (slight precision change)
$ws.Range("F574").Value = 45929.47446090278
$ws.Range("G574").Value = 45929.47446027778

# Append new rows 575-596
$ws.Range("A575").Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Range("B575").Value = "MUTAG"
$ws.Range("C575").Value = "feature-KNN"
$ws.Range("D575").Value = "feature-KNN_trained_on_MUTAG.joblib"
$ws.Range("F575").Value = 45934.69965002315
$ws.Range("G575").Value = 45934.6996496875
$ws.Range("F575:G575").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A576").Value = "Fucntionality_test_MUTAG_with_SVC_Hybrid-Prototype-GED_poly"
$ws.Range("B576").Value = "MUTAG"
$ws.Range("C576").Value = "SVC_Hybrid-Prototype-GED_poly"
$ws.Range("D576").Value = "SVC_Hybrid-Prototype-GED_poly_trained_on_MUTAG.joblib"
$ws.Range("F576").Value = 45934.73216295139
$ws.Range("G576").Value = 45934.73216259259
$ws.Range("F576:G576").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A577").Value = "Fucntionality_test_PTC_FR_with_SVC_Hybrid-Prototype-GED_poly"
$ws.Range("B577").Value = "PTC_FR"
$ws.Range("C577").Value = "SVC_Hybrid-Prototype-GED_poly"
$ws.Range("D577").Value = "SVC_Hybrid-Prototype-GED_poly_trained_on_PTC_FR.joblib"
$ws.Range("F577").Value = 45934.74649184028
$ws.Range("G577").Value = 45934.74649076389
$ws.Range("F577:G577").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A578").Value = "Fucntionality_test_MUTAG_with_SVC_VertexHistogram_rbf"
$ws.Range("B578").Value = "MUTAG"
$ws.Range("C578").Value = "SVC_VertexHistogram_rbf"
$ws.Range("D578").Value = "SVC_VertexHistogram_rbf_trained_on_MUTAG.joblib"
$ws.Range("F578").Value = 45934.77812075231
$ws.Range("G578").Value = 45934.77806451389
$ws.Range("F578:G578").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A579").Value = "Fucntionality_test_MUTAG_with_SVC_EdgeHistogram_rbf"
$ws.Range("B579").Value = "MUTAG"
$ws.Range("C579").Value = "SVC_EdgeHistogram_rbf"
$ws.Range("D579").Value = "SVC_EdgeHistogram_rbf_trained_on_MUTAG.joblib"
$ws.Range("F579").Value = 45934.78180104167
$ws.Range("G579").Value = 45934.78180091435
$ws.Range("F579:G579").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A580").Value = "Fucntionality_test_MUTAG_with_SVC_EdgeHistogram_rbf"
$ws.Range("B580").Value = "MUTAG"
$ws.Range("C580").Value = "SVC_EdgeHistogram_rbf"
$ws.Range("D580").Value = "SVC_EdgeHistogram_rbf_trained_on_MUTAG.joblib"
$ws.Range("F580").Value = 45934.78180104167
$ws.Range("G580").Value = 45934.78180091435
$ws.Range("F580:G580").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A581").Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Range("B581").Value = "MUTAG"
$ws.Range("C581").Value = "SVC_CombinedHistogram_rbf"
$ws.Range("D581").Value = "SVC_CombinedHistogram_rbf_trained_on_MUTAG.joblib"
$ws.Range("F581").Value = 45934.84572197917
$ws.Range("G581").Value = 45934.8457218287
$ws.Range("F581:G581").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A582").Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Range("B582").Value = "MUTAG"
$ws.Range("C582").Value = "SVC_CombinedHistogram_rbf"
$ws.Range("D582").Value = "SVC_CombinedHistogram_rbf_trained_on_MUTAG.joblib"
$ws.Range("F582").Value = 45934.84572197917
$ws.Range("G582").Value = 45934.8457218287
$ws.Range("F582:G582").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A583").Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Range("B583").Value = "MUTAG"
$ws.Range("C583").Value = "SVC_CombinedHistogram_rbf"
$ws.Range("D583").Value = "SVC_CombinedHistogram_rbf_trained_on_MUTAG.joblib"
$ws.Range("F583").Value = 45934.85872717592
$ws.Range("G583").Value = 45934.85872704861
$ws.Range("F583:G583").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A584").Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Range("B584").Value = "MUTAG"
$ws.Range("C584").Value = "SVC_CombinedHistogram_rbf"
$ws.Range("D584").Value = "SVC_CombinedHistogram_rbf_trained_on_MUTAG.joblib"
$ws.Range("F584").Value = 45934.85872717592
$ws.Range("G584").Value = 45934.85872704861
$ws.Range("F584:G584").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A585").Value = "Fucntionality_test_MUTAG_with_SVC_CombinedHistogram_rbf"
$ws.Range("B585").Value = "MUTAG"
$ws.Range("C585").Value = "SVC_CombinedHistogram_rbf"
$ws.Range("D585").Value = "SVC_CombinedHistogram_rbf_trained_on_MUTAG.joblib"
$ws.Range("F585").Value = 45934.86840295139
$ws.Range("G585").Value = 45934.86824746527
$ws.Range("F585:G585").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A586").Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Range("B586").Value = "MUTAG"
$ws.Range("C586").Value = "feature-KNN"
$ws.Range("D586").Value = "feature-KNN_trained_on_MUTAG.joblib"
$ws.Range("F586").Value = 45934.87197950231
$ws.Range("G586").Value = 45934.87197908565
$ws.Range("F586:G586").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A587").Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Range("B587").Value = "MUTAG"
$ws.Range("C587").Value = "feature-KNN"
$ws.Range("D587").Value = "feature-KNN_trained_on_MUTAG.joblib"
$ws.Range("F587").Value = 45934.87640118055
$ws.Range("G587").Value = 45934.87640077547
$ws.Range("F587:G587").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A588").Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Range("B588").Value = "MUTAG"
$ws.Range("C588").Value = "feature-KNN"
$ws.Range("D588").Value = "feature-KNN_trained_on_MUTAG.joblib"
$ws.Range("F588").Value = 45934.88721842592
$ws.Range("G588").Value = 45934.88689296296
$ws.Range("F588:G588").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A589").Value = "Fucntionality_test_MUTAG_with_feature-KNN"
$ws.Range("B589").Value = "MUTAG"
$ws.Range("C589").Value = "feature-KNN"
$ws.Range("D589").Value = "feature-KNN_trained_on_MUTAG.joblib"
$ws.Range("F589").Value = 45934.88872850694
$ws.Range("G589").Value = 45934.88872775463
$ws.Range("F589:G589").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A590").Value = "Fucntionality_test_MUTAG_with_SVC_Hybrid-Prototype-GED_poly"
$ws.Range("B590").Value = "MUTAG"
$ws.Range("C590").Value = "SVC_Hybrid-Prototype-GED_poly"
$ws.Range("D590").Value = "SVC_Hybrid-Prototype-GED_poly_trained_on_MUTAG.joblib"
$ws.Range("F590").Value = 45934.98598952546
$ws.Range("G590").Value = 45934.98598822916
$ws.Range("F590:G590").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A591").Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Range("B591").Value = "MUTAG"
$ws.Range("C591").Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Range("D591").Value = "SVC_Random-Walk-Edit_precomputed_trained_on_MUTAG.joblib"
$ws.Range("F591").Value = 45935.12126365741
$ws.Range("G591").Value = 45935.12125907408
$ws.Range("F591:G591").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A592").Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Range("B592").Value = "MUTAG"
$ws.Range("C592").Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Range("D592").Value = "SVC_Random-Walk-Edit_precomputed_trained_on_MUTAG.joblib"
$ws.Range("F592").Value = 45935.12517533565
$ws.Range("G592").Value = 45935.12517086806
$ws.Range("F592:G592").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A593").Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Range("B593").Value = "MUTAG"
$ws.Range("C593").Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Range("D593").Value = "SVC_Random-Walk-Edit_precomputed_trained_on_MUTAG.joblib"
$ws.Range("F593").Value = 45935.12592638889
$ws.Range("G593").Value = 45935.12592527777
$ws.Range("F593:G593").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A594").Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Range("B594").Value = "MUTAG"
$ws.Range("C594").Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Range("D594").Value = "SVC_Random-Walk-Edit_precomputed_trained_on_MUTAG.joblib"
$ws.Range("F594").Value = 45935.12736097222
$ws.Range("G594").Value = 45935.12736019676
$ws.Range("F594:G594").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A595").Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Range("B595").Value = "MUTAG"
$ws.Range("C595").Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Range("D595").Value = "SVC_Random-Walk-Edit_precomputed_trained_on_MUTAG.joblib"
$ws.Range("F595").Value = 45935.70941407407
$ws.Range("G595").Value = 45935.70941309028
$ws.Range("F595:G595").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A596").Value = "Fucntionality_test_MUTAG_with_SVC_Random-Walk-Edit_precomputed"
$ws.Range("B596").Value = "MUTAG"
$ws.Range("C596").Value = "SVC_Random-Walk-Edit_precomputed"
$ws.Range("D596").Value = "SVC_Random-Walk-Edit_precomputed_trained_on_MUTAG.joblib"
$ws.Range("F596").Value = 45935.70991768294
$ws.Range("G596").Value = 45935.70991684647
$ws.Range("F596:G596").NumberFormat = "YYYY-MM-DD HH:MM:SS"
